$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "29.263.54"
$ws.Range("E2").Value = "  +0.69%  "
$ws.Range("D3").Value = "1.859.60"
$ws.Range("E3").Value = "  +0.89%  "
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "0.9997"
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = "  -0.05%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "0.7030"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +0.27%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "238.00"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  +0.27%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "1.000"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  -0.04%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.08086"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  +8.67%  "
$ws.Range("E9").Value = "  +0.43%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "23.33"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  +0.37%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.08185"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  +0.96%  "
$ws.Range("D12").Value = "1.855.19"
$ws.Range("E12").Value = "  +1.89%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "0.7166"
$ws.Range("D13").Style = "Normal"
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "5.188"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  -1.00%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "89.28"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  +0.42%  "
$ws.Range("D16").Value = "29.276.45"
$ws.Range("E16").Value = "  +1.14%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "5.775"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  -0.12%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "13.38"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  +2.96%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "0.000007833"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  +2.41%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "236.17"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  -1.43%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "0.9992"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  -0.06%  "
$ws.Range("D22").Value = "2.105.97"
$ws.Range("E22").Value = "  +2.29%  "
$ws.Range("E23").Value = "  -0.03%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "7.467"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  -0.91%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "161.93"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  +0.50%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "8.983"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  +0.60%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "0.1457"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  -0.14%  "
$ws.Range("E28").Value = "  +0.75%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "1.995"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  +3.26%  "
$ws.Range("E30").Value = "  +4.51%  "
$ws.Range("B31").Value = "Filecoin"
$ws.Range("C31").Value = "https://coinranking.com/coin/ymQub4fuB+filecoin-fil"
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "4.409"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  -1.41%  "
$ws.Range("B32").Value = "PancakeSwap"
$ws.Range("C32").Value = "https://coinranking.com/coin/ncYFcP709+pancakeswap-cake"
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "1.484"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  -0.46%  "
$ws.Range("E33").Value = "  +1.19%  "
$ws.Range("E34").Value = "  +1.10%  "
$ws.Range("E35").Value = "  -0.78%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.7087"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  +0.39%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "1.003"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  -1.38%  "
$ws.Range("E38").Value = "  +0.57%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.01849"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  -0.90%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "2.719"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  +1.83%  "
$ws.Range("D41").Value = "1.146.25"
$ws.Range("E41").Value = "  +8.28%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.9224"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  +3.22%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "5.964"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  +0.10%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.4279"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  +0.09%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "70.95"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  +1.76%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.9994"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  -0.06%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "103.33"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  +1.87%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "1.780"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  +1.81%  "
$ws.Range("D49").Value = "2.003.14"
$ws.Range("E49").Value = "  +1.88%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "9.196"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  +0.63%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "6.938"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  -1.41%  "
